$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 0) Capture the OLD "_GoBack" bookmark position first (before any
#    other edits shift character offsets around).
# -----------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBmPos = $oldBm.Start

# -----------------------------------------------------------------
# 1) "intro@rijnijssel.nl " -> "DORP-Logistics " (with auto color)
# -----------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$targetPara = $null
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*intro@rijnijssel.nl*") {
        $targetPara = $p
        break
    }
}

$full = $targetPara.Range.Text
$idx = $full.IndexOf("intro@rijnijssel.nl ")
$oldWord = "intro@rijnijssel.nl "
$newWord = "DORP-Logistics "
$absStart = $targetPara.Range.Start + $idx
$absEnd = $absStart + $oldWord.Length

$introRange = $d.Range($absStart, $absEnd)
$introRange.Text = $newWord
$introRange.Font.Color = -16777216

# Track how much the document shifted because of the replace above,
# so we can translate the old bookmark's captured position.
$delta = $newWord.Length - $oldWord.Length
if ($absStart -lt $oldBmPos) {
    $oldBmPos = $oldBmPos + $delta
}

# Add the new "_GoBack" bookmark right after the replaced text
# (matching real Word's behaviour of moving "_GoBack" to the most
# recent edit location).
$bmRange = $d.Range($introRange.End, $introRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# -----------------------------------------------------------------
# 2) Remove the stray "_GoBack" bookmark that used to sit inside the
#    version-history table ("21-02-2018" cell).
# -----------------------------------------------------------------
$delRange = $d.Range($oldBmPos - 1, $oldBmPos + 1)
$savedText = $delRange.Text
$delRange.Text = ""

$insRange = $d.Range($oldBmPos - 1, $oldBmPos - 1)
$insRange.InsertBefore($savedText)
